$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.930.44"
Set-TextValue $ws.Range("E2") "  +0.12%  "

Set-TextValue $ws.Range("D3") "1.876.97"
Set-TextValue $ws.Range("E3") "  -0.65%  "

Set-TextValue $ws.Range("E4") "  +0.06%  "

Set-TextValue $ws.Range("D5") "0.7432"
Set-TextValue $ws.Range("E5") "  -4.11%  "

Set-TextValue $ws.Range("D6") "242.54"
Set-TextValue $ws.Range("E6") "  -0.17%  "

Set-TextValue $ws.Range("D7") "1.001"
Set-TextValue $ws.Range("E7") "  +0.05%  "

Set-TextValue $ws.Range("D8") "0.3149"
Set-TextValue $ws.Range("E8") "  +1.07%  "

Set-TextValue $ws.Range("D9") "0.07216"
Set-TextValue $ws.Range("E9") "  +0.59%  "

Set-TextValue $ws.Range("D10") "24.67"
Set-TextValue $ws.Range("E10") "  -3.76%  "

Set-TextValue $ws.Range("D11") "0.08389"
Set-TextValue $ws.Range("E11") "  -2.51%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.921.39"
Set-TextValue $ws.Range("E12") "  -1.95%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D13") "0.7525"
Set-TextValue $ws.Range("E13") "  -1.45%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "5.416"
Set-TextValue $ws.Range("E14") "  +0.70%  "

Set-TextValue $ws.Range("D15") "92.66"
Set-TextValue $ws.Range("E15") "  -1.20%  "

Set-TextValue $ws.Range("D16") "29.946.93"
Set-TextValue $ws.Range("E16") "  +0.00%  "

Set-TextValue $ws.Range("D17") "6.076"
Set-TextValue $ws.Range("E17") "  -1.65%  "

Set-TextValue $ws.Range("D18") "253.24"
Set-TextValue $ws.Range("E18") "  +3.56%  "

Set-TextValue $ws.Range("D19") "13.59"
Set-TextValue $ws.Range("E19") "  -1.36%  "

Set-TextValue $ws.Range("D20") "0.000007855"
Set-TextValue $ws.Range("E20") "  +0.43%  "

Set-TextValue $ws.Range("D21") "1.000"
Set-TextValue $ws.Range("E21") "  +0.18%  "

Set-TextValue $ws.Range("D22") "2.127.84"
Set-TextValue $ws.Range("E22") "  -2.90%  "

Set-TextValue $ws.Range("D23") "8.039"
Set-TextValue $ws.Range("E23") "  +0.46%  "

Set-TextValue $ws.Range("E24") "  +0.00%  "

Set-TextValue $ws.Range("D25") "0.1559"
Set-TextValue $ws.Range("E25") "  -4.96%  "

Set-TextValue $ws.Range("D26") "9.267"
Set-TextValue $ws.Range("E26") "  -1.24%  "

Set-TextValue $ws.Range("D27") "165.09"
Set-TextValue $ws.Range("E27") "  +1.88%  "

Set-TextValue $ws.Range("D28") "18.76"
Set-TextValue $ws.Range("E28") "  -0.06%  "

Set-TextValue $ws.Range("D29") "2.035"
Set-TextValue $ws.Range("E29") "  -0.35%  "

Set-TextValue $ws.Range("D30") "1.518"
Set-TextValue $ws.Range("E30") "  +5.13%  "

Set-TextValue $ws.Range("D31") "4.605"
Set-TextValue $ws.Range("E31") "  +1.82%  "

Set-TextValue $ws.Range("D32") "1.534"
Set-TextValue $ws.Range("E32") "  -0.42%  "

Set-TextValue $ws.Range("D33") "4.281"
Set-TextValue $ws.Range("E33") "  +4.38%  "

Set-TextValue $ws.Range("D34") "0.05329"
Set-TextValue $ws.Range("E34") "  -1.93%  "

Set-TextValue $ws.Range("D35") "1.236"
Set-TextValue $ws.Range("E35") "  -0.30%  "

Set-TextValue $ws.Range("D36") "0.7499"
Set-TextValue $ws.Range("E36") "  +0.51%  "

Set-TextValue $ws.Range("D37") "0.9998"
Set-TextValue $ws.Range("E37") "  -0.41%  "

Set-TextValue $ws.Range("D38") "2.691"
Set-TextValue $ws.Range("E38") "  -0.16%  "

Set-TextValue $ws.Range("D39") "0.01965"
Set-TextValue $ws.Range("E39") "  -0.01%  "

Set-TextValue $ws.Range("D40") "2.756"
Set-TextValue $ws.Range("E40") "  -0.98%  "

Set-TextValue $ws.Range("D41") "0.4538"
Set-TextValue $ws.Range("E41") "  +1.49%  "

Set-TextValue $ws.Range("D42") "1.112.47"
Set-TextValue $ws.Range("E42") "  -0.35%  "

Set-TextValue $ws.Range("D43") "6.082"
Set-TextValue $ws.Range("E43") "  -0.09%  "

Set-TextValue $ws.Range("E44") "  -0.94%  "

Set-TextValue $ws.Range("D45") "0.8569"
Set-TextValue $ws.Range("E45") "  +0.66%  "

Set-TextValue $ws.Range("D46") "1.002"
Set-TextValue $ws.Range("E46") "  +0.17%  "

Set-TextValue $ws.Range("D47") "103.66"
Set-TextValue $ws.Range("E47") "  +1.15%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D48") "1.856"
Set-TextValue $ws.Range("E48") "  -0.56%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D49") "7.617"
Set-TextValue $ws.Range("E49") "  -0.35%  "

Set-TextValue $ws.Range("D50") "2.025.58"
Set-TextValue $ws.Range("E50") "  -2.94%  "

Set-TextValue $ws.Range("D51") "2.898"
Set-TextValue $ws.Range("E51") "  -2.82%  "
